$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename / retarget columns for Jang's framework validation data ---
$ws.Range("A1").Value = "P_in (Pa)"
$ws.Range("B1").Value = "volume flow rate (m^3/s)"
$ws.Range("C1").Value = "Pressure in (Pa)"
$ws.Range("D1").Value = "Pressure out (Pa)"
$ws.Range("E1").Value = "delP "

# --- Replace the validation data table (rows 2-32) ---
$data = New-Object 'object[,]' 31,4
$data[0,0] = [double]"100"
$data[0,1] = [double]"7.0768e-08"
$data[0,2] = [double]"99.832"
$data[0,3] = [double]"0.084902"
$data[1,0] = [double]"125.89"
$data[1,1] = [double]"9.8258e-08"
$data[1,2] = [double]"125.68"
$data[1,3] = [double]"0.10779"
$data[2,0] = [double]"158.49"
$data[2,1] = [double]"1.3655e-07"
$data[2,2] = [double]"158.22"
$data[2,3] = [double]"0.13545"
$data[3,0] = [double]"199.53"
$data[3,1] = [double]"1.8988e-07"
$data[3,2] = [double]"199.18"
$data[3,3] = [double]"0.16992"
$data[4,0] = [double]"251.19"
$data[4,1] = [double]"2.6396e-07"
$data[4,2] = [double]"250.75"
$data[4,3] = [double]"0.21286"
$data[5,0] = [double]"316.23"
$data[5,1] = [double]"3.6674e-07"
$data[5,2] = [double]"315.68"
$data[5,3] = [double]"0.26827"
$data[6,0] = [double]"398.11"
$data[6,1] = [double]"5.0938e-07"
$data[6,2] = [double]"397.42"
$data[6,3] = [double]"0.33866"
$data[7,0] = [double]"501.19"
$data[7,1] = [double]"7.0786e-07"
$data[7,2] = [double]"500.32"
$data[7,3] = [double]"0.42628"
$data[8,0] = [double]"630.96"
$data[8,1] = [double]"9.8342e-07"
$data[8,2] = [double]"629.87"
$data[8,3] = [double]"0.53782"
$data[9,0] = [double]"794.33"
$data[9,1] = [double]"1.3659e-06"
$data[9,2] = [double]"792.97"
$data[9,3] = [double]"0.67841"
$data[10,0] = [double]"1000"
$data[10,1] = [double]"1.8986e-06"
$data[10,2] = [double]"998.28"
$data[10,3] = [double]"0.85088"
$data[11,0] = [double]"1258.9"
$data[11,1] = [double]"2.6375e-06"
$data[11,2] = [double]"1256.8"
$data[11,3] = [double]"1.0738"
$data[12,0] = [double]"1584.9"
$data[12,1] = [double]"3.6675e-06"
$data[12,2] = [double]"1582.2"
$data[12,3] = [double]"1.3409"
$data[13,0] = [double]"1995.3"
$data[13,1] = [double]"5.0917e-06"
$data[13,2] = [double]"1991.9"
$data[13,3] = [double]"1.7035"
$data[14,0] = [double]"2511.9"
$data[14,1] = [double]"7.0757e-06"
$data[14,2] = [double]"2507.6"
$data[14,3] = [double]"2.1427"
$data[15,0] = [double]"3162.3"
$data[15,1] = [double]"9.8354e-06"
$data[15,2] = [double]"3156.8"
$data[15,3] = [double]"2.6886"
$data[16,0] = [double]"3981.1"
$data[16,1] = [double]"1.3668e-05"
$data[16,2] = [double]"3974.2"
$data[16,3] = [double]"3.3805"
$data[17,0] = [double]"5011.9"
$data[17,1] = [double]"1.8992e-05"
$data[17,2] = [double]"5003.3"
$data[17,3] = [double]"4.2578"
$data[18,0] = [double]"6309.6"
$data[18,1] = [double]"2.6401e-05"
$data[18,2] = [double]"6298.7"
$data[18,3] = [double]"5.3395"
$data[19,0] = [double]"7943.3"
$data[19,1] = [double]"3.6739e-05"
$data[19,2] = [double]"7929.5"
$data[19,3] = [double]"6.611"
$data[20,0] = [double]"10000"
$data[20,1] = [double]"5.1052e-05"
$data[20,2] = [double]"9982.7"
$data[20,3] = [double]"8.3503"
$data[21,0] = [double]"12589"
$data[21,1] = [double]"7.1101e-05"
$data[21,2] = [double]"12567"
$data[21,3] = [double]"10.263"
$data[22,0] = [double]"15849"
$data[22,1] = [double]"9.9011e-05"
$data[22,2] = [double]"15821"
$data[22,3] = [double]"12.669"
$data[23,0] = [double]"19953"
$data[23,1] = [double]"0.00013808"
$data[23,2] = [double]"19917"
$data[23,3] = [double]"15.442"
$data[24,0] = [double]"25119"
$data[24,1] = [double]"0.00019289"
$data[24,2] = [double]"25073"
$data[24,3] = [double]"18.474"
$data[25,0] = [double]"31623"
$data[25,1] = [double]"0.00027024"
$data[25,2] = [double]"31563"
$data[25,3] = [double]"21.43"
$data[26,0] = [double]"39811"
$data[26,1] = [double]"0.00038052"
$data[26,2] = [double]"39731"
$data[26,3] = [double]"23.394"
$data[27,0] = [double]"50119"
$data[27,1] = [double]"0.00053993"
$data[27,2] = [double]"50012"
$data[27,3] = [double]"22.769"
$data[28,0] = [double]"63096"
$data[28,1] = [double]"0.00077638"
$data[28,2] = [double]"62946"
$data[28,3] = [double]"16.433"
$data[29,0] = [double]"79433"
$data[29,1] = [double]"0.0011527"
$data[29,2] = [double]"79181"
$data[29,3] = [double]"-3.4836"
$data[30,0] = [double]"100000"
$data[30,1] = [double]"0.0018028"
$data[30,2] = [double]"99543"
$data[30,3] = [double]"-35.56"
$ws.Range("A2:D32").Value = $data

# --- delP formula (D-C) still spans the shrunk data range ---
$ws.Range("E2").Formula = "=D2-C2"
$ws.Range("E3:E32").Formula = "=D3-C3"

# --- Scientific notation for the now very-small flow-rate figures ---
$ws.Range("B16:B30").NumberFormat = "0.00E+00"
$ws.Range("A32").NumberFormat = "0.00E+00"

# --- The table now ends at row 32; blank out the old trailing rows 33-35 ---
$ws.Range("A33:B35").ClearContents()
$ws.Range("C33:E35").ClearContents()

# --- Match the author's final selection ---
$ws.Range("E2:E32").Select()
